$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Remove the old rows 7-14 (their jobs are no longer in the feed) ---
$ws.Rows("7:14").Delete()

# --- Resize columns B (48 -> 44) and H (20 -> 21) ---
# ColumnWidth is expressed in "characters"; the engine adds a constant
# 5/6-character padding when it stores the OOXML <col width="..."/>, so we
# subtract it here to land exactly on the target stored width.
$ws.Columns("B").ColumnWidth = 44 - 5/6
$ws.Columns("H").ColumnWidth = 21 - 5/6

# --- Row 2: new job posting ---
$ws.Range("A2").Value = "2025-11-09 06:24:36"
$ws.Range("B2").Value = "【急募】Wordpressを用いた比較サイトの新規開発"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5430121"
$ws.Range("G2").Value = 123
$ws.Range("H2").Value = "◆開発 ◇サイト ○WordPress"

# --- Row 3: new job posting ---
$ws.Range("A3").Value = "2025-11-09 06:24:36"
$ws.Range("B3").Value = "Glideメインで作成したシステムをLinux+MySQL型に移行するための新規開発"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5430095"
$ws.Range("G3").Value = 115
$ws.Range("H3").Value = "◆開発 ◇MySQL"

# --- Row 4: new job posting (no skill tags) ---
$ws.Range("A4").Value = "2025-11-09 06:24:36"
$ws.Range("B4").Value = "インターネット情報収集(selenium)"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5430171"
$ws.Range("G4").Value = 10
$ws.Range("H4").ClearContents()

# --- Row 5: new job posting (no skill tags) ---
$ws.Range("A5").Value = "2025-11-09 06:24:36"
$ws.Range("B5").Value = "MT4 RSXを使用したEAの作成依頼"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5430008"
$ws.Range("G5").Value = 10
$ws.Range("H5").ClearContents()

# --- Row 6: new job posting (no skill tags) ---
$ws.Range("A6").Value = "2025-11-09 06:24:36"
$ws.Range("B6").Value = "【急募】LINE × QRコード連携で自動取得設定を実現!"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5430015"
$ws.Range("G6").Value = 10
$ws.Range("H6").ClearContents()

# --- Hyperlinks: rebuild F2:F6 against their new URLs ---
# The old collection still references the removed rows/old URLs, so clear
# it out entirely and re-add just the 5 links that remain.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5430121")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5430095")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5430171")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5430008")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5430015")

# Hyperlinks.Add() re-applies the "Hyperlink" cell style via a freshly
# minted xf record; reassert the shared built-in "Hyperlink" style so the
# cells keep referencing the workbook's existing style (matches original).
$ws.Range("F2:F6").Style = "Hyperlink"
